# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stale forecast values that should no longer be populated
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Correct small floating point precision drift in remaining forecast values
$ws.Range("C4").Value = 6.277541464866965
$ws.Range("E4").Value = 9.926356894614985

$ws.Range("C5").Value = 6.535114773304795

$ws.Range("C11").Value = 2.508469427909921
$ws.Range("E11").Value = 2.540874511056646

$ws.Range("E12").Value = 3.025650759930021

$ws.Range("E13").Value = 0.676128192849057

$ws.Range("C14").Value = 3.047037961814514

$ws.Range("C15").Value = -0.22288476972816

$ws.Range("C16").Value = -1.165854108406639

$ws.Range("C19").Value = 2.039329803030099
